$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the new dates to be stored as plain text (shared strings) rather
# than being auto-converted to date serial numbers by Excel's input parser.
$ws.Range("A10:A12").NumberFormat = "@"

$ws.Range("A10").Value = "2024-10-04"
$ws.Range("B10").Value = 0.01804

$ws.Range("A11").Value = "2024-10-05"
$ws.Range("B11").Value = 0.01886

$ws.Range("A12").Value = "2024-03-09"
$ws.Range("B12").Value = 0.01732

# Restore the default (unstyled) formatting on the new cells so the written
# cells stay style-free, matching the rest of the data rows.
$ws.Range("A10:A12").Style = $ws.Range("A2").Style
